$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 values
$a3 = $ws.Range("A3").Value2
$b3 = $ws.Range("B3").Value2
$d3 = $ws.Range("D3").Value2
$e3 = $ws.Range("E3").Value2
$f3 = $ws.Range("F3").Value2
$g3 = $ws.Range("G3").Value2
$h3 = $ws.Range("H3").Value2

# Row 4 values
$a4 = $ws.Range("A4").Value2
$b4 = $ws.Range("B4").Value2
$d4 = $ws.Range("D4").Value2
$e4 = $ws.Range("E4").Value2
$f4 = $ws.Range("F4").Value2
$g4 = $ws.Range("G4").Value2
$h4 = $ws.Range("H4").Value2

# Swap row 3 and row 4 for columns A, B, D, E, F, G, H
$ws.Range("A3").Value2 = $a4
$ws.Range("B3").Value2 = $b4
$ws.Range("D3").Value2 = $d4
$ws.Range("E3").Value2 = $e4
$ws.Range("F3").Value2 = $f4
$ws.Range("G3").Value2 = $g4
$ws.Range("H3").Value2 = $h4

$ws.Range("A4").Value2 = $a3
$ws.Range("B4").Value2 = $b3
$ws.Range("D4").Value2 = $d3
$ws.Range("E4").Value2 = $e3
$ws.Range("F4").Value2 = $f3
$ws.Range("G4").Value2 = $g3
$ws.Range("H4").Value2 = $h3

# Row 5: only B5 changes
$ws.Range("B5").Value2 = 90805
